# Updates the cryptos price list with fresh market data (prices + 1h volume %)
# fetched by the GitHub Actions scraper. Two coins (NEARProtocol and
# Binance-PegBSC-USD) also swapped rank positions (rows 31/32), so those two
# rows get their Coin name, Link, Price and Volume(1h) cells rewritten too.
#
# Price values that are plain numbers (e.g. "6.55") are written with a
# leading apostrophe so Excel stores them as literal text instead of
# reinterpreting them as numbers -- matching how the sheet already stores
# every Price/Volume cell as text (note some prices, like "69.385.47", use
# '.' as a thousands separator and are not valid numbers at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.385.47"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.678.09"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'686.22"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'159.28"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.494"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "'7.06"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").Value = "'0.436"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "'0.0000233"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "4.294.49"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'32.33"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").Value = "3.668.24"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "69.370.02"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'15.84"
$ws.Range("E18").Value = "  -3.05%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'469.88"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").Value = "'9.96"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'0.650"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "'79.81"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "3.821.51"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D26").Value = "'0.0000124"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("E28").Value = "  -4.02%  "
$ws.Range("D29").Value = "'2.71"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").Value = "'1.74"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.55"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "  -5.89%  "
$ws.Range("D34").Value = "'26.90"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").Value = "3.649.98"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").Value = "'8.17"
$ws.Range("E37").Value = "  -4.67%  "
$ws.Range("D38").Value = "'6.12"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D40").Value = "'2.22"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'0.0899"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'166.56"
$ws.Range("E43").Value = "  +5.52%  "
$ws.Range("D44").Value = "'0.940"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D46").Value = "'0.000285"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "'2.72"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("D49").Value = "'1.30"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "'27.36"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'7.79"
$ws.Range("E51").Value = "  -3.87%  "
